$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto market data (price + volume %) scraped on Tue Oct 17 2023.
# Cells in column D occasionally look like plain numbers (e.g. "24.12"); Excel
# auto-converts such text to a Number on assignment, which would both change the
# cell type and round/trim the display (e.g. "1.00" -> 1). Forcing a Text number
# format before the assignment keeps these as literal strings, matching the source.

$ws.Range('D2').Value = '28.385.69'
$ws.Range('E2').Value = '  +4.22%  '
$ws.Range('D3').Value = '1.596.10'
$ws.Range('E3').Value = '  +1.96%  '
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.72'
$ws.Range('E5').Value = '  +1.87%  '
$ws.Range('E6').Value = '  +1.13%  '
$ws.Range('E7').Value = '  -0.17%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '24.12'
$ws.Range('E8').Value = '  +8.90%  '
$ws.Range('E9').Value = '  +0.65%  '
$ws.Range('E10').Value = '  +1.00%  '
$ws.Range('E11').Value = '  +1.86%  '
$ws.Range('D12').Value = '1.825.09'
$ws.Range('E12').Value = '  +2.01%  '
$ws.Range('D13').Value = '1.583.89'
$ws.Range('E13').Value = '  +1.49%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.532'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.75'
$ws.Range('E15').Value = '  +0.04%  '
$ws.Range('D16').Value = '28.403.46'
$ws.Range('E16').Value = '  +4.49%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.23'
$ws.Range('E17').Value = '  +2.18%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '228.31'
$ws.Range('E18').Value = '  +4.78%  '
$ws.Range('D19').Value = '0.0₃0711'
$ws.Range('E19').Value = '  +1.35%  '
$ws.Range('E20').Value = '  +0.61%  '
$ws.Range('E21').Value = '  -0.19%  '
$ws.Range('E22').Value = '  -0.77%  '
$ws.Range('E23').Value = '  -0.17%  '
$ws.Range('E24').Value = '  +0.68%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '151.67'
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.22'
$ws.Range('E26').Value = '  +1.33%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.107'
$ws.Range('E27').Value = '  +0.68%  '
$ws.Range('E29').Value = '  -0.13%  '
$ws.Range('E30').Value = '  +0.80%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0476'
$ws.Range('E31').Value = '  +1.38%  '
$ws.Range('E32').Value = '  +0.05%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.16'
$ws.Range('E33').Value = '  -0.27%  '
$ws.Range('D34').Value = '1.399.90'
$ws.Range('E34').Value = '  -4.10%  '
$ws.Range('E35').Value = '  -1.53%  '
$ws.Range('E36').Value = '  -4.99%  '
$ws.Range('E37').Value = '  +0.22%  '
$ws.Range('E38').Value = '  +0.70%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.54'
$ws.Range('E39').Value = '  +8.35%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.542'
$ws.Range('E40').Value = '  +0.58%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.816'
$ws.Range('E41').Value = '  +0.35%  '
$ws.Range('E42').Value = '  -2.71%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.91'
$ws.Range('E43').Value = '  +8.94%  '
$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.00'
$ws.Range('E44').Value = '  -0.22%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.985'
$ws.Range('E45').Value = '  +0.11%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '64.51'
$ws.Range('E46').Value = '  +0.12%  '
$ws.Range('D47').Value = '1.735.97'
$ws.Range('E47').Value = '  +2.08%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '87.63'
$ws.Range('E48').Value = '  +2.09%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.14'
$ws.Range('E49').Value = '  +0.03%  '
$ws.Range('E50').Value = '  -1.12%  '
$ws.Range('E51').Value = '  +0.18%  '
